# Optimizacion de nomenclaturas numericas
#
# The protocol-summary table on the worksheet had its numeric message codes
# renumbered/simplified (8->4, 9->5, 10->6, 11->7, 12->4) and a couple of
# mislabeled/duplicated rows were corrected: a "Solicitar lista de
# conectados" row label was restored in A5 and a "Salir de partida" label
# was added to the previously-blank A9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Desconexión) / Row 4 (Iniciar sesión): notification column now
# references the "3$nombre1/.../nombreN" code instead of the old
# duplicated/incorrect text.
$ws.Range("D2").Value = '3$nombre1/…/nombreN'
$ws.Range("D4").Value = '3$nombre1/…/nombreN'

# Row 5: restore the proper operation label.
$ws.Range("A5").Value = 'Solicitar lista de conectados'

# Row 6 (Invitar jugadores): renumber 8 -> 4 / 8 -> 4 / 12 -> 5
$ws.Range("A6").Value = 'Invitar jugadores'
$ws.Range("B6").Value = '4/N/nombre1/…/nombreN'
$ws.Range("C6").Value = '4$ID_partida'
$ws.Range("D6").Value = '5$nombre_host/ID_partida'

# Row 7 (Respuesta a invitación): renumber 9 -> 5 / 10 -> 6
$ws.Range("A7").Value = 'Respuesta a invitación'
$ws.Range("B7").Value = '5/ID_partida/respuesta'
$ws.Range("D7").Value = '6$ID_partida/mensaje'

# Row 8 (Mensaje en chat): renumber 10 -> 6
$ws.Range("A8").Value = 'Mensaje en chat'
$ws.Range("B8").Value = '6/ID_partida/mensaje'
$ws.Range("D8").Value = '6$ID_partida/mensaje'

# Row 9 (Salir de partida): add missing label, renumber 11 -> 7
$ws.Range("A9").Value = 'Salir de partida'
$ws.Range("B9").Value = '7/ID_partida'

# Restore the active selection to B1, matching the saved worksheet view.
$ws.Range("B1").Select()
